# Journal_progress.docx - "Commit efter tabelnr - Test"
#
# The author added table numbers/captions to the two result tables
# ("Tabel 1: ..." / "Tabel 2: ...") and centered those caption
# paragraphs. (The rest of the upstream diff - style-id spelling,
# added xmlns:mo/mv namespaces and w:proofErr bookmarks - is inert
# save-environment noise from the editing Word install and carries no
# visible content change, so it is not reproduced here.)

$d = $word.ActiveDocument

# --- Table 1 caption: "Suggestion 1 (Mixed realization)" ----------------
# becomes "Tabel 1: Suggestion 1 (Mixed realization)", centered.
$rng = $d.Content
$rng.Find.Execute("Suggestion 1 (Mixed realization)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "Tabel 1: Suggestion 1 (Mixed realization)"
$rng.Paragraphs.Item(1).Alignment = 1

# --- Table 2 caption: "Suggestion 2 (SW realization)" --------------------
# gets "Tabel 2: " prefixed (as a separate leading run), centered.
$rng2 = $d.Content
$rng2.Find.Execute("Suggestion 2 (SW realization)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.InsertBefore("Tabel 2: ")

$rng3 = $d.Content
$rng3.Find.Execute("Tabel 2: Suggestion 2 (SW realization)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Paragraphs.Item(1).Alignment = 1

# --- Drop the stale lastRenderedPageBreak pagination-cache markers -------
# (these four paragraphs are exactly the ones the upstream diff shows
# losing their <w:lastRenderedPageBreak/> marker; re-running them through
# Find/Replace rebuilds the run and drops the stale marker without
# altering any visible text).
$d.Content.Find.Execute("Configure individual libraries for block realization in HW.", $true, $false, $false, $false, $false, $true, 1, $false, "Configure individual libraries for block realization in HW.", 2)
$d.Content.Find.Execute("If the remote user changes the volume twice", $true, $false, $false, $false, $false, $true, 1, $false, "If the remote user changes the volume twice", 2)
$d.Content.Find.Execute("SW interface for updating bass, treble and volume coefficients.", $true, $false, $false, $false, $false, $true, 1, $false, "SW interface for updating bass, treble and volume coefficients.", 2)
$d.Content.Find.Execute("15/12-2010", $true, $false, $false, $false, $false, $true, 1, $false, "15/12-2010", 2)

Write-Output "edit complete"
